$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.997677963090024
$ws.Range("D2").Value = 0.00003789130113614438
$ws.Range("E2").Value = 331.2468016463672
$ws.Range("F2").Value = 1447.583428519429
$ws.Range("G2").Value = 1116.336626873062
$ws.Range("H2").Value = 46708.45839110907
$ws.Range("I2").Value = 4597.466859655824
$ws.Range("J2").Value = 119.5719141440429
$ws.Range("K2").Value = 11.72069730627138
$ws.Range("L2").Value = 0.00003789130113614438
$ws.Range("M2").Value = 0.4113311823267408
$ws.Range("O2").Value = 0.001961902034617951
$ws.Range("P2").Value = 0.4068401668521797
$ws.Range("Q2").Value = 0.4154143515655463
$ws.Range("S2").Value = 0.01239709188923505
$ws.Range("T2").Value = 0.3629135192468673
$ws.Range("U2").Value = 0.9976347172963363
$ws.Range("V2").Value = 0.9977212126331396
$ws.Range("W2").Value = 131.2926114503143
$ws.Range("C3").Value = 0.9976721691882303
$ws.Range("D3").Value = 0.00003507449044370555
$ws.Range("E3").Value = 331.248725335829
$ws.Range("F3").Value = 1447.60280615436
$ws.Range("G3").Value = 1116.354080818531
$ws.Range("H3").Value = 46656.49410210855
$ws.Range("I3").Value = 4618.306479344215
$ws.Range("J3").Value = 126.6308857560833
$ws.Range("K3").Value = 12.55986337145032
$ws.Range("L3").Value = 0.00003507449044370555
$ws.Range("M3").Value = 0.4134911095642975
$ws.Range("O3").Value = 0.002044577747540558
$ws.Range("P3").Value = 0.4154159704865861
$ws.Range("Q3").Value = 0.4228271923107911
$ws.Range("S3").Value = 0.01143706524846398
$ws.Range("T3").Value = 0.3459185720153358
$ws.Range("U3").Value = 0.9976315661259776
$ws.Range("V3").Value = 0.9977127755556626
$ws.Range("W3").Value = 139.1907491275336
$ws.Range("C4").Value = 0.997660382840662
$ws.Range("D4").Value = 0.00003625817038997172
$ws.Range("E4").Value = 331.2526387042212
$ws.Range("F4").Value = 1447.624676270078
$ws.Range("G4").Value = 1116.372037565857
$ws.Range("H4").Value = 44966.09125519858
$ws.Range("I4").Value = 4501.07228842146
$ws.Range("J4").Value = 125.0162487015378
$ws.Range("K4").Value = 13.11176800510184
$ws.Range("L4").Value = 0.00003625817038997172
$ws.Range("M4").Value = 0.4137356048528412
$ws.Range("O4").Value = 0.002148399943571912
$ws.Range("P4").Value = 0.4223033224530864
$ws.Range("Q4").Value = 0.4174042203422526
$ws.Range("S4").Value = 0.01181690488845993
$ws.Range("T4").Value = 0.3466455613965865
$ws.Range("U4").Value = 0.997618324174933
$ws.Range("V4").Value = 0.9977024450528494
$ws.Range("W4").Value = 138.1280167066397
$ws.Range("C5").Value = 0.9976371141494025
$ws.Range("D5").Value = 0.00003081009375132641
$ws.Range("E5").Value = 331.2603647754245
$ws.Range("F5").Value = 1447.661843630309
$ws.Range("G5").Value = 1116.401478854885
$ws.Range("H5").Value = 45456.77684533728
$ws.Range("I5").Value = 4766.617979168981
$ws.Range("J5").Value = 120.2870410475835
$ws.Range("K5").Value = 13.53173955983582
$ws.Range("L5").Value = 0.00003081009375132641
$ws.Range("M5").Value = 0.4188400336602033
$ws.Range("O5").Value = 0.00207296107132909
$ws.Range("P5").Value = 0.4060749588531597
$ws.Range("Q5").Value = 0.4178537556066292
$ws.Range("S5").Value = 0.009993427511747251
$ws.Range("T5").Value = 0.3538787377604689
$ws.Range("U5").Value = 0.9976007758491043
$ws.Range("V5").Value = 0.9976734550970926
$ws.Range("W5").Value = 133.8187806074194
$ws.Range("C6").Value = 0.9976267020778394
$ws.Range("D6").Value = 0.00002735640832469095
$ws.Range("E6").Value = 331.2638220872796
$ws.Range("F6").Value = 1447.68142542972
$ws.Range("G6").Value = 1116.41760334244
$ws.Range("H6").Value = 44852.64703610058
$ws.Range("I6").Value = 4826.049369025101
$ws.Range("J6").Value = 118.6625581292928
$ws.Range("K6").Value = 15.14051744672575
$ws.Range("L6").Value = 0.00002735640832469094
$ws.Range("M6").Value = 0.4185683107223698
$ws.Range("O6").Value = 0.002086874027805218
$ws.Range("P6").Value = 0.4054751808071732
$ws.Range("Q6").Value = 0.4207737288794502
$ws.Range("S6").Value = 0.008818628863412906
$ws.Range("T6").Value = 0.3419238068643196
$ws.Range("U6").Value = 0.9975938603909323
$ws.Range("V6").Value = 0.9976595459271734
$ws.Range("W6").Value = 133.8030755760185
$ws.Range("C7").Value = 0.9976186448894745
$ws.Range("D7").Value = 0.00002640795616512301
$ws.Range("E7").Value = 331.2664975134323
$ws.Range("F7").Value = 1447.700564549888
$ws.Range("G7").Value = 1116.434067036456
$ws.Range("H7").Value = 44710.71958289205
$ws.Range("I7").Value = 4641.207546625555
$ws.Range("J7").Value = 119.085833272801
$ws.Range("K7").Value = 13.95164970143268
$ws.Range("L7").Value = 0.00002640795616512301
$ws.Range("M7").Value = 0.4193689835665058
$ws.Range("O7").Value = 0.002112949713188917
$ws.Range("P7").Value = 0.4069321082466584
$ws.Range("Q7").Value = 0.4310228552161054
$ws.Range("S7").Value = 0.008489063097475182
$ws.Range("T7").Value = 0.305319665793264
$ws.Range("U7").Value = 0.9975867176428636
$ws.Range("V7").Value = 0.9976505741797809
$ws.Range("W7").Value = 133.0374829742337
$ws.Range("C8").Value = 0.9976037012348412
$ws.Range("D8").Value = 0.00002619841071991989
$ws.Range("E8").Value = 331.2714597365318
$ws.Range("F8").Value = 1447.727525712897
$ws.Range("G8").Value = 1116.456065976366
$ws.Range("H8").Value = 45580.57290947381
$ws.Range("I8").Value = 4807.138936609352
$ws.Range("J8").Value = 120.6288937242428
$ws.Range("K8").Value = 13.25113063849876
$ws.Range("L8").Value = 0.00002619841071991989
$ws.Range("M8").Value = 0.4202625593096239
$ws.Range("O8").Value = 0.002176955563387554
$ws.Range("P8").Value = 0.399731134980214
$ws.Range("Q8").Value = 0.4370924994946145
$ws.Range("S8").Value = 0.008401320537696648
$ws.Range("T8").Value = 0.280698927497818
$ws.Range("U8").Value = 0.9975718464217815
$ws.Range("V8").Value = 0.997635558082364
$ws.Range("W8").Value = 133.8800243627416
$ws.Range("C9").Value = 0.9975832338571962
$ws.Range("D9").Value = 0.00002931386485278933
$ws.Range("E9").Value = 331.2782564205972
$ws.Range("F9").Value = 1447.755332045013
$ws.Range("G9").Value = 1116.477075624416
$ws.Range("H9").Value = 46890.16841601417
$ws.Range("I9").Value = 4402.013476330178
$ws.Range("J9").Value = 124.0622908978741
$ws.Range("K9").Value = 14.63063901418081
$ws.Range("L9").Value = 0.00002931386485278933
$ws.Range("M9").Value = 0.4232655949531891
$ws.Range("O9").Value = 0.002149523804893702
$ws.Range("P9").Value = 0.3970859917496964
$ws.Range("Q9").Value = 0.4439784987947744
$ws.Range("S9").Value = 0.009470161696208255
$ws.Range("T9").Value = 0.2379996211213599
$ws.Range("U9").Value = 0.9975482445494183
$ws.Range("V9").Value = 0.9976182256195812
$ws.Range("W9").Value = 138.692929912055
$ws.Range("C10").Value = 0.9975757172569771
$ws.Range("D10").Value = 0.00003242554057759752
$ws.Range("E10").Value = 331.2807525581552
$ws.Range("F10").Value = 1447.769638079903
$ws.Range("G10").Value = 1116.488885521748
$ws.Range("H10").Value = 49191.44110915018
$ws.Range("I10").Value = 4322.430490424518
$ws.Range("J10").Value = 129.5949813472934
$ws.Range("K10").Value = 15.08669094085558
$ws.Range("L10").Value = 0.00003242554057759752
$ws.Range("M10").Value = 0.4229639667734412
$ws.Range("O10").Value = 0.002108440016721441
$ws.Range("P10").Value = 0.3932426057399453
$ws.Range("Q10").Value = 0.4543394629788884
$ws.Range("S10").Value = 0.01053300200787734
$ws.Range("T10").Value = 0.2280451654286512
$ws.Range("U10").Value = 0.9975376519172418
$ws.Range("V10").Value = 0.9976137855019169
$ws.Range("W10").Value = 144.681672288149
$ws.Range("C11").Value = 0.9975698417594873
$ws.Range("D11").Value = 0.0000316410132567224
$ws.Range("E11").Value = 331.2827037390637
$ws.Range("F11").Value = 1447.782064648375
$ws.Range("G11").Value = 1116.499360909312
$ws.Range("H11").Value = 48819.08900488963
$ws.Range("I11").Value = 4613.979699295842
$ws.Range("J11").Value = 122.976872173722
$ws.Range("K11").Value = 15.16246746501848
$ws.Range("L11").Value = 0.0000316410132567224
$ws.Range("M11").Value = 0.4232857045363292
$ws.Range("O11").Value = 0.002058472253117562
$ws.Range("P11").Value = 0.3929510585828793
$ws.Range("Q11").Value = 0.4517106241735679
$ws.Range("S11").Value = 0.01027801247799601
$ws.Range("T11").Value = 0.2524598166617361
$ws.Range("U11").Value = 0.997532695096764
$ws.Range("V11").Value = 0.9976069911888886
$ws.Range("W11").Value = 138.1393396387405
